$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# --- Rename "Exchange" test case to "Replace" (row 6) ---
$ws.Range("B6").Value = "Replace a product"
$ws.Range("C6").Value = "Replace an existing product for another one (+ secret key)"

# --- Fill in the previously-blank row 11 (Nr. 10) with new test case content ---
$ws.Range("B11").Value = "decrease amount of Item"
$ws.Range("C11").Value = "decrease an amount of the chosen Item in the Vending Machine"
$ws.Range("D11").Value = "The Amount will decrease from the Vending machine"

# --- Add 3 new rows to Table1 for the new test cases ---
$row12 = $lo.ListRows.Add()
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Stop the User from buying"
$ws.Range("C12").Value = "If the Vending Machine is Empty, the User won't able to so anything"
$ws.Range("D12").Value = "At the buying Process, if empty, won't able to do anything"
$ws.Rows.Item(12).RowHeight = 22.5

$row13 = $lo.ListRows.Add()
$ws.Range("A11").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "able to quit buying any Time"
$ws.Range("C13").Value = "The User shoulb be able to quit in the Process Stage"
$ws.Range("D13").Value = "inputing -1 should stop the process and return the Money"
$ws.Rows.Item(13).RowHeight = 22.5

$row14 = $lo.ListRows.Add()
$ws.Range("A11").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 13

$excel.CutCopyMode = 0

# --- Update the active selection to match the saved view ---
$ws.Range("F4").Select()
